# "analisi del progetto" correction
# - Row 3 (item "2", "aggiungere un test") and Row 4 (item "3", "cancellare un
#   test"): TIPOLOGIA cell text "funzionale" stays the same, but the
#   spell-check proofErr wrapper around it is dropped by re-writing the cell
#   text.
# - Row 8 (item "7", "Esportare in formato CSV..."): TIPOLOGIA "tecnologico" -> "funzionale"
# - Row 9 (item "8", "Salvare i dati su un file binario"): TIPOLOGIA "tecnologico" -> "funzionale"
# - Row 10 (item "9", "Caricare i dati all'avvio..."): TIPOLOGIA "tecnologico" -> "Non funzionale"

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

$table.Cell(3, 3).Range.Text = "funzionale"
$table.Cell(4, 3).Range.Text = "funzionale"
$table.Cell(8, 3).Range.Text = "funzionale"
$table.Cell(9, 3).Range.Text = "funzionale"
$table.Cell(10, 3).Range.Text = "Non funzionale"
